# Update "想去人数" (interest count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5162
$wsExhibit.Range("F6").Value = 301
$wsExhibit.Range("F7").Value = 784
$wsExhibit.Range("F8").Value = 267

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5162
$wsAll.Range("F6").Value = 301
$wsAll.Range("F7").Value = 784
$wsAll.Range("F9").Value = 267
